# StructureDefinition-match-compared-to.xlsx: re-point the IG at the new
# "linuxforhealth.org" home, bump the version/date/publisher metadata, and
# drop the stray duplicated FHIR constraint text that used to sit on the
# top-level "Extension" row of the Elements sheet.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-compared-to"  # URL
$wsMeta.Range("B3").Value = "8.0.0"                                                                      # Version
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"                                                  # Date
$wsMeta.Range("B9").Value = "LinuxForHealth Team"                                                        # Publisher

# --- Elements sheet ---------------------------------------------------------
$wsElements = $wb.Worksheets.Item("Elements")

# "Extension.url" (row 3) has a Fixed Value cell that repeats the same
# canonical URL string as Metadata!B2 - keep it consistent with the new URL.
$wsElements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-compared-to"

# The top-level "Extension" row (row 2) had the ele-1/ext-1 constraint text
# duplicated into its Constraint(s) column; clear it so it only appears on
# the rows that actually declare it (e.g. Extension.extension, row 4).
$wsElements.Range("AI2").Value = ""
